$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep text (string) storage, matching the inline-string cells
# in the original workbook (e.g. '0.999', '1.00' must not become numeric).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.015.09'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.306.79'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.28'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.95'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.502'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.20'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +7.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.54'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +12.37%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.117'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.94'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.661.68'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.298.46'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.802'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.905.34'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.44'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +6.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.24'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.99'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.44'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.26'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +14.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.46'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.83'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.70'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.69'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.19'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.03'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.81'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.68'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.85%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0694'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.85'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.79'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.110'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.993.09'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0290'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.26%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.23'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.71%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.92'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.85%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.75'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '56.16'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +6.21%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.55'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.12%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.526.26'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.52'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.07%  '
